# Generate Report for Handoff
# Adds a new row (for file ffff8723fde2-cb08-45b2-b7b0-653a0738b867.md) to each
# of the three sheets (Overview, zh-cn, de-de), and refreshes the handoff
# bookkeeping fields (new GUID / timestamps / xlf hash) for the existing row.

$wb = $excel.ActiveWorkbook

$oldGuid = "b9b62266-5f5a-4ea8-a245-52d179bcbf78"
$newGuid = "60400108-2d7c-4903-8f3e-1d13301ecf94"
$newFileGuid = "ffff8723fde2-cb08-45b2-b7b0-653a0738b867"

$oldHash = "acf5a99546d2ece5fa2a38219468f29fdbf9c500"
$newHash = "05d5dfb294c4aea78dc6e4d4c1ca3ff93948ee14"

$newHandoffDateTime = "2016-03-20 17:25:40"
$newXlfDateTime = "2016-03-20 17:25:32"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Update existing row with regenerated identifiers.
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("D2").Value = $newHandoffDateTime

# Append the new file's row.
$ws.Range("A3").Value = "$newFileGuid.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $newHandoffDateTime

$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("D2").Style

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f574c4c07b27bc992281868a26ff3b8dce9ae2c4/e2e/$newFileGuid.md", "", "", "$newFileGuid.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("E2").Value = $newXlfDateTime

# Append the new file's row (same handoff metadata as the refreshed row above).
$ws.Range("A3").Value = "$newFileGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("E3").Value = $newXlfDateTime
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("D2").Style
$ws.Range("E3").Style = $ws.Range("E2").Style
$ws.Range("H3").Style = $ws.Range("H2").Style
$ws.Range("J3").Style = $ws.Range("J2").Style

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f574c4c07b27bc992281868a26ff3b8dce9ae2c4/e2e/$newFileGuid.md", "", "", "$newFileGuid.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5988c4f95ceb06797bf9f0a43469047efb64ecf7/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$newGuid.$newHash.zh-cn.xlf", "", "", "$newGuid.$newHash.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"

# Append the new file's row (same handoff metadata as the refreshed row above).
$ws.Range("A3").Value = "$newFileGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("E3").Value = $ws.Range("E2").Value
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("D2").Style
$ws.Range("E3").Style = $ws.Range("E2").Style
$ws.Range("H3").Style = $ws.Range("H2").Style
$ws.Range("J3").Style = $ws.Range("J2").Style

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f574c4c07b27bc992281868a26ff3b8dce9ae2c4/e2e/$newFileGuid.md", "", "", "$newFileGuid.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43b808e55a146c46ae84bf3819d652f953488d64/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$newGuid.$newHash.de-de.xlf", "", "", "$newGuid.$newHash.de-de.xlf")
